$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in FILENAME / MODEL columns for rows 3-16 (same values as row 2),
# matching the "model answers" added alongside the updated evaluation scores.
for ($r = 3; $r -le 16; $r++) {
    $ws.Cells.Item($r, 1).Value = "Nicholls-Diver-Finding"
    $ws.Cells.Item($r, 2).Value = "gemma3"
}

# Leave the selection where the user finished typing (column B, the MODEL
# column that was just filled in down to the last row).
[void]$ws.Range("B2:B16").Select()
